$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 83 (R2 / WDS_ST_EXISTING) entirely, shifting all rows below it up by one.
$ws.Rows(83).Delete()
